$wb = $excel.ActiveWorkbook

# Rename the old "Sheet3" (4th sheet) to "vacancies"
$vacancies = $wb.Worksheets.Item("Sheet3")
$vacancies.Name = "vacancies"

# Add a new worksheet after "vacancies" (i.e. at the end) and name it "Employee"
$employee = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$employee.Name = "Employee"

# Fill in the Employee sheet with the new data
$employee.Range("A1").Value = "First name"
$employee.Range("B1").Value = "Middle name"
$employee.Range("C1").Value = "Last name"
$employee.Range("D1").Value = "username"

$employee.Range("A2").Value = "Abu"
$employee.Range("B2").Value = "Mohd"
$employee.Range("C2").Value = "Khan"
$employee.Range("D2").Value = "abumadkhan"

$employee.Range("A3").Value = "Ahmad"
$employee.Range("B3").Value = "Mohd"
$employee.Range("C3").Value = "Khan"
$employee.Range("D3").Value = "ahmadkhan"

# Match the saved selection on the new sheet
$null = $employee.Range("J8").Select()
